$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.48051789516684
$ws.Range("D2").Value = 8.371820453742533
$ws.Range("E2").Value = 27.1116297535703
$ws.Range("F2").Value = 46.89940864320836
$ws.Range("G2").Value = 3.667457791298832
$ws.Range("L2").Value = 11.71159954251972
$ws.Range("M2").Value = 18.07308932590204
$ws.Range("B3").Value = 19.2342359464687
$ws.Range("D3").Value = 8.321892989255966
$ws.Range("E3").Value = 25.55785238654527
$ws.Range("F3").Value = 45.86389900863831
$ws.Range("G3").Value = 3.67617734045005
$ws.Range("L3").Value = 11.51769649563437
$ws.Range("M3").Value = 17.92457234033711
$ws.Range("B4").Value = 19.08666690820657
$ws.Range("D4").Value = 8.295644523553619
$ws.Range("E4").Value = 24.55223641622176
$ws.Range("F4").Value = 45.24262715832337
$ws.Range("G4").Value = 3.681779998004051
$ws.Range("L4").Value = 11.40021730890678
$ws.Range("M4").Value = 17.83722668216217
$ws.Range("B5").Value = 19.02751445644548
$ws.Range("D5").Value = 8.286047885664685
$ws.Range("E5").Value = 24.12960347405355
$ws.Range("F5").Value = 44.99345751862929
$ws.Range("G5").Value = 3.684126138855605
$ws.Range("L5").Value = 11.35279867377712
$ws.Range("M5").Value = 17.80262841496755
$ws.Range("B6").Value = 19.01775349289236
$ws.Range("D6").Value = 8.284520502494679
$ws.Range("E6").Value = 24.05865506880773
$ws.Range("F6").Value = 44.95233491875253
$ws.Range("G6").Value = 3.684519532491883
$ws.Range("L6").Value = 11.34495402662495
$ws.Range("M6").Value = 17.79694435855617
$ws.Range("B7").Value = 19.08586509169778
$ws.Range("D7").Value = 8.295510659751043
$ws.Range("E7").Value = 24.54658840587652
$ws.Range("F7").Value = 45.2392501041525
$ws.Range("G7").Value = 3.681811383157545
$ws.Range("L7").Value = 11.39957588525881
$ws.Range("M7").Value = 17.83675600871624
$ws.Range("B8").Value = 19.39488237739787
$ws.Range("D8").Value = 8.353686478413051
$ws.Range("E8").Value = 26.58666542016152
$ws.Range("F8").Value = 46.53955949990327
$ws.Range("G8").Value = 3.670412919236955
$ws.Range("L8").Value = 11.64444927417863
$ws.Range("M8").Value = 18.02110074359832
$ws.Range("B9").Value = 20.02686358177789
$ws.Range("D9").Value = 8.502987991414573
$ws.Range("E9").Value = 30.17367589717518
$ws.Range("F9").Value = 49.1892396624972
$ws.Range("G9").Value = 3.650014360827535
$ws.Range("L9").Value = 12.13460710445433
$ws.Range("M9").Value = 18.4117494322087
$ws.Range("B10").Value = 20.50301217798077
$ws.Range("D10").Value = 8.634316422876861
$ws.Range("E10").Value = 32.55327785132847
$ws.Range("F10").Value = 51.17608034056897
$ws.Range("G10").Value = 3.636189661085651
$ws.Range("L10").Value = 12.49741517452095
$ws.Range("M10").Value = 18.71462435593434
$ws.Range("B11").Value = 20.72133740082091
$ws.Range("D11").Value = 8.698743139515543
$ws.Range("E11").Value = 33.57995677902237
$ws.Range("F11").Value = 52.08442348445992
$ws.Range("G11").Value = 3.630146423916453
$ws.Range("L11").Value = 12.66235708832335
$ws.Range("M11").Value = 18.85542788416213
$ws.Range("B12").Value = 20.80418841715914
$ws.Range("D12").Value = 8.723808504256857
$ws.Range("E12").Value = 33.96069954985811
$ws.Range("F12").Value = 52.42870447471233
$ws.Range("G12").Value = 3.627892821621618
$ws.Range("L12").Value = 12.72474637293058
$ws.Range("M12").Value = 18.90914438707884
$ws.Range("B13").Value = 20.78633826845364
$ws.Range("D13").Value = 8.718380610518606
$ws.Range("E13").Value = 33.87905734165221
$ws.Range("F13").Value = 52.35454876547622
$ws.Range("G13").Value = 3.628376633544483
$ws.Range("L13").Value = 12.71131373322744
$ws.Range("M13").Value = 18.89755851109278
$ws.Range("B14").Value = 20.72815048268974
$ws.Range("D14").Value = 8.700791924942221
$ws.Range("E14").Value = 33.61144183786438
$ws.Range("F14").Value = 52.11274332883544
$ws.Range("G14").Value = 3.629960323281109
$ws.Range("L14").Value = 12.66749162075203
$ws.Range("M14").Value = 18.85983942814319
$ws.Range("B15").Value = 20.69252959727261
$ws.Range("D15").Value = 8.690105194673583
$ws.Range("E15").Value = 33.44647265833795
$ws.Range("F15").Value = 51.9646612913855
$ws.Range("G15").Value = 3.630934901556547
$ws.Range("L15").Value = 12.64063848857629
$ws.Range("M15").Value = 18.83678601371424
$ws.Range("B16").Value = 20.48877261577917
$ws.Range("D16").Value = 8.630199850987509
$ws.Range("E16").Value = 32.48505749049276
$ws.Range("F16").Value = 51.11677956525328
$ws.Range("G16").Value = 3.636589503102231
$ws.Range("L16").Value = 12.48662894785714
$ws.Range("M16").Value = 18.70548016231989
$ws.Range("B17").Value = 20.36416411790059
$ws.Range("D17").Value = 8.594646661679782
$ws.Range("E17").Value = 31.88095201536541
$ws.Range("F17").Value = 50.59754263242407
$ws.Range("G17").Value = 3.640120999291869
$ws.Range("L17").Value = 12.39208251580396
$ws.Range("M17").Value = 18.62567621538975
$ws.Range("B18").Value = 20.29265899364379
$ws.Range("D18").Value = 8.574638819656819
$ws.Range("E18").Value = 31.52823283107745
$ws.Range("F18").Value = 50.29933591697647
$ws.Range("G18").Value = 3.642175377596887
$ws.Range("L18").Value = 12.33769559213155
$ws.Range("M18").Value = 18.58006203268981
$ws.Range("B19").Value = 20.26847937172811
$ws.Range("D19").Value = 8.567940447024895
$ws.Range("E19").Value = 31.40790644735333
$ws.Range("F19").Value = 50.19845521929516
$ws.Range("G19").Value = 3.642874946464544
$ws.Range("L19").Value = 12.31928177571648
$ws.Range("M19").Value = 18.56466827542605
$ws.Range("B20").Value = 20.37741220868882
$ws.Range("D20").Value = 8.598385706465558
$ws.Range("E20").Value = 31.94580394317965
$ws.Range("F20").Value = 50.65277278218166
$ws.Range("G20").Value = 3.639742672545044
$ws.Range("L20").Value = 12.40214821226334
$ws.Range("M20").Value = 18.6341420562018
$ws.Range("B21").Value = 20.74523743430164
$ws.Range("D21").Value = 8.705940061732949
$ws.Range("E21").Value = 33.69026518716237
$ws.Range("F21").Value = 52.18376161247288
$ws.Range("G21").Value = 3.629494213596526
$ws.Range("L21").Value = 12.68036558958858
$ws.Range("M21").Value = 18.87090795692928
$ws.Range("B22").Value = 20.98662610325524
$ws.Range("D22").Value = 8.780124887686853
$ws.Range("E22").Value = 34.78352148578163
$ws.Range("F22").Value = 53.18602108382894
$ws.Range("G22").Value = 3.622999096678653
$ws.Range("L22").Value = 12.86176072577585
$ws.Range("M22").Value = 19.02794608037353
$ws.Range("B23").Value = 20.85772483247966
$ws.Range("D23").Value = 8.740177233485626
$ws.Range("E23").Value = 34.20431755045924
$ws.Range("F23").Value = 52.65105069258339
$ws.Range("G23").Value = 3.626447266585689
$ws.Range("L23").Value = 12.76500451766466
$ws.Range("M23").Value = 18.943934119604
$ws.Range("B24").Value = 20.37142232412719
$ws.Range("D24").Value = 8.596693937857903
$ws.Range("E24").Value = 31.9165012392876
$ws.Range("F24").Value = 50.62780224149947
$ws.Range("G24").Value = 3.639913639044776
$ws.Range("L24").Value = 12.3975976029239
$ws.Range("M24").Value = 18.6303138176638
$ws.Range("B25").Value = 19.85352579044729
$ws.Range("D25").Value = 8.458800465124556
$ws.Range("E25").Value = 29.24803874209861
$ws.Range("F25").Value = 48.4639083421573
$ws.Range("G25").Value = 3.655326512464335
$ws.Range("L25").Value = 12.00129330054927
$ws.Range("M25").Value = 18.30313995596493
